# Apply cryptos list update (crypto price/volume refresh).
# Source diff changes 90 cell values across rows 2-51 (data rows),
# including two name/link/price swaps (rows 19<->20 and 50<->51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.637.43'
$ws.Cells.Item(3, 4).Value = '1.598.05'
$ws.Cells.Item(3, 5).Value = '  +0.34%  '
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$ws.Cells.Item(5, 4).Value = '''211.48'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -0.30%  '
$ws.Cells.Item(6, 5).Value = '  +0.35%  '
$ws.Cells.Item(7, 5).Value = '  +0.07%  '
$ws.Cells.Item(8, 4).Value = '''0.0618'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = '  +0.11%  '
$ws.Cells.Item(9, 5).Value = '  -0.11%  '
$ws.Cells.Item(10, 4).Value = '''19.46'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -0.72%  '
$ws.Cells.Item(11, 5).Value = '  +0.27%  '
$ws.Cells.Item(12, 4).Value = '1.822.14'
$ws.Cells.Item(13, 4).Value = '1.605.41'
$ws.Cells.Item(13, 5).Value = '  +0.89%  '
$ws.Cells.Item(14, 5).Value = '  -0.21%  '
$ws.Cells.Item(15, 5).Value = '  -0.28%  '
$ws.Cells.Item(16, 4).Value = '''64.79'
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = '  -0.60%  '
$ws.Cells.Item(17, 4).Value = '26.634.58'
$ws.Cells.Item(17, 5).Value = '  +0.09%  '
$ws.Cells.Item(18, 4).Value = '0.0₃0735'
$ws.Cells.Item(18, 5).Value = '  +0.81%  '
$ws.Cells.Item(19, 2).Value = 'Dai'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(19, 4).Value = '''1.00'
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = '  +0.07%  '
$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).Value = '''208.80'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  +0.12%  '
$ws.Cells.Item(21, 4).Value = '''7.06'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +5.51%  '
$ws.Cells.Item(22, 5).Value = '  -0.09%  '
$ws.Cells.Item(23, 5).Value = '  -0.89%  '
$ws.Cells.Item(24, 4).Value = '''8.88'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  -0.03%  '
$ws.Cells.Item(25, 4).Value = '''145.58'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -0.24%  '
$ws.Cells.Item(26, 5).Value = '  +0.12%  '
$ws.Cells.Item(27, 4).Value = '''7.16'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +0.16%  '
$ws.Cells.Item(28, 5).Value = '  +0.72%  '
$ws.Cells.Item(29, 4).Value = '''15.24'
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = '  -0.46%  '
$ws.Cells.Item(30, 4).Value = '''0.0507'
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +0.54%  '
$ws.Cells.Item(31, 4).Value = '''1.15'
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -0.20%  '
$ws.Cells.Item(32, 5).Value = '  -0.31%  '
$ws.Cells.Item(33, 4).Value = '''2.93'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +0.94%  '
$ws.Cells.Item(34, 4).Value = '''0.625'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -6.29%  '
$ws.Cells.Item(35, 4).Value = '1.271.79'
$ws.Cells.Item(35, 5).Value = '  -2.12%  '
$ws.Cells.Item(36, 5).Value = '  +0.35%  '
$ws.Cells.Item(37, 5).Value = '  +0.05%  '
$ws.Cells.Item(38, 5).Value = '  -0.74%  '
$ws.Cells.Item(39, 4).Value = '''0.841'
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +1.40%  '
$ws.Cells.Item(40, 4).Value = '''5.49'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +2.54%  '
$ws.Cells.Item(41, 5).Value = '  +1.08%  '
$ws.Cells.Item(42, 4).Value = '''0.787'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -0.66%  '
$ws.Cells.Item(43, 4).Value = '''64.09'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +1.45%  '
$ws.Cells.Item(44, 4).Value = '''0.945'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +18.06%  '
$ws.Cells.Item(45, 4).Value = '1.734.82'
$ws.Cells.Item(45, 5).Value = '  +0.31%  '
$ws.Cells.Item(46, 4).Value = '''89.99'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = '  +1.06%  '
$ws.Cells.Item(47, 4).Value = '''1.60'
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = '  +0.10%  '
$ws.Cells.Item(48, 4).Value = '''0.102'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +4.28%  '
$ws.Cells.Item(49, 4).Value = '''0.0508'
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = '  +1.12%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = '''7.48'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -0.48%  '
$ws.Cells.Item(51, 2).Value = 'USDD'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Cells.Item(51, 4).Value = '''1.00'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +0.23%  '
